$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PV-Test-03")

# Rename header cells to reflect new column naming scheme
$ws.Range("A1").Value = "Row ID"
$ws.Range("C1").Value = "Task"
$ws.Range("E1").Value = "Start Date"
$ws.Range("F1").Value = "End Date"

# Update the active selection on the sheet
$ws.Activate()
$ws.Range("E2").Select()
